# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to H..N columns across several rows in multiple sheets
# of the Ixion_Profits workbook, per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3080564.5
$ws.Range("I17").Value = 163
$ws.Range("J17").Value = 3135571.8
$ws.Range("K17").Value = 489
$ws.Range("L17").Value = 9406715.399999999
$ws.Range("M17").Value = -321
$ws.Range("N17").Value = -9407051.399999999
# Row 113
$ws.Range("H113").Value = 3364.5
$ws.Range("I113").Value = 3442
$ws.Range("J113").Value = 3235.3333
$ws.Range("K113").Value = 3442
$ws.Range("L113").Value = 3235.3333
$ws.Range("M113").Value = -188
$ws.Range("N113").Value = -9743.3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4447.8413
$ws.Range("I32").Value = 2550.7163
$ws.Range("J32").Value = 12921.667
$ws.Range("K32").Value = 2550.7163
$ws.Range("L32").Value = 12921.667
$ws.Range("M32").Value = -2263.7163
# Row 45
$ws.Range("H45").Value = 6842.1763
$ws.Range("I45").Value = 8080.4287
$ws.Range("J45").Value = 1063.6666
$ws.Range("K45").Value = 8080.4287
$ws.Range("L45").Value = 1063.6666
$ws.Range("M45").Value = -7703.4287
# Row 61
$ws.Range("H61").Value = 1669.7188
$ws.Range("I61").Value = 1361.3572
$ws.Range("J61").Value = 3828.25
$ws.Range("K61").Value = 1361.3572
$ws.Range("L61").Value = 3828.25
$ws.Range("M61").Value = -1149.3572
$ws.Range("N61").Value = -4252.25
# Row 74
$ws.Range("H74").Value = 1159.8823
$ws.Range("I74").Value = 1054.8572
$ws.Range("J74").Value = 1650
$ws.Range("K74").Value = 1054.8572
$ws.Range("L74").Value = 1650
$ws.Range("M74").Value = -180.8571999999999
$ws.Range("N74").Value = -3398
# Row 77
$ws.Range("H77").Value = 1159.8823
$ws.Range("I77").Value = 1054.8572
$ws.Range("J77").Value = 1650
$ws.Range("K77").Value = 5274.286
$ws.Range("L77").Value = 8250
$ws.Range("M77").Value = -906.2860000000001
$ws.Range("N77").Value = -16986
# Row 122
$ws.Range("H122").Value = 1070686.8
$ws.Range("I122").Value = 1605334.2
$ws.Range("J122").Value = 1391.625
$ws.Range("K122").Value = 4816002.6
$ws.Range("L122").Value = 4174.875
$ws.Range("M122").Value = -4813552.6
$ws.Range("N122").Value = -9074.875
# Row 132
$ws.Range("H132").Value = 3209.1516
$ws.Range("I132").Value = 1996.3478
$ws.Range("J132").Value = 5998.6
$ws.Range("K132").Value = 5989.0434
$ws.Range("L132").Value = 17995.8
$ws.Range("M132").Value = -3459.0434
$ws.Range("N132").Value = -23055.8
# Row 136
$ws.Range("H136").Value = 1669.7188
$ws.Range("I136").Value = 1361.3572
$ws.Range("J136").Value = 3828.25
$ws.Range("K136").Value = 4084.0716
$ws.Range("L136").Value = 11484.75
$ws.Range("M136").Value = -1534.0716
$ws.Range("N136").Value = -16584.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3357.8372
$ws.Range("I31").Value = 1748.9259
$ws.Range("J31").Value = 6072.875
$ws.Range("K31").Value = 1748.9259
$ws.Range("L31").Value = 6072.875
$ws.Range("M31").Value = -1453.9259
$ws.Range("N31").Value = -6662.875
# Row 34
$ws.Range("H34").Value = 3357.8372
$ws.Range("I34").Value = 1748.9259
$ws.Range("J34").Value = 6072.875
$ws.Range("K34").Value = 1748.9259
$ws.Range("L34").Value = 6072.875
$ws.Range("M34").Value = -1546.9259
$ws.Range("N34").Value = -6476.875
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
# Row 99
$ws.Range("H99").Value = 5126
$ws.Range("I99").Value = 3349.2856
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 3349.2856
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -1851.2856
$ws.Range("N99").Value = -32996
# Row 126
$ws.Range("H126").Value = 5126
$ws.Range("I126").Value = 3349.2856
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 10047.8568
$ws.Range("L126").Value = 90000
$ws.Range("M126").Value = -7577.856800000001
$ws.Range("N126").Value = -94940
# Row 132
$ws.Range("H132").Value = 2007.35
$ws.Range("I132").Value = 1753.8857
$ws.Range("J132").Value = 3781.6
$ws.Range("K132").Value = 5261.6571
$ws.Range("L132").Value = 11344.8
$ws.Range("M132").Value = -2731.6571
$ws.Range("N132").Value = -16404.8
# Row 134
$ws.Range("H134").Value = 1139.44
$ws.Range("I134").Value = 1116.6511
$ws.Range("J134").Value = 1279.4286
$ws.Range("K134").Value = 3349.9533
$ws.Range("L134").Value = 3838.2858
$ws.Range("M134").Value = -814.9533000000001
$ws.Range("N134").Value = -8908.2858

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 3846244.5
$ws.Range("I12").Value = 11111216
$ws.Range("J12").Value = 82.82353000000001
$ws.Range("K12").Value = 33333648
$ws.Range("L12").Value = 248.47059
$ws.Range("M12").Value = -33333475
$ws.Range("N12").Value = -594.47059
# Row 99
$ws.Range("H99").Value = 725
$ws.Range("I99").Value = 725
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2175
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 71
# Row 103
$ws.Range("H103").Value = 4343.4116
$ws.Range("I103").Value = 8068.2
$ws.Range("J103").Value = 2791.4167
$ws.Range("K103").Value = 24204.6
$ws.Range("L103").Value = 8374.250100000001
$ws.Range("M103").Value = -23325.6
$ws.Range("N103").Value = -10132.2501

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 31987.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 31987.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 31987.5
$ws.Range("N64").Value = -32483.5
# Row 67
$ws.Range("H67").Value = 31987.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 31987.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 31987.5
$ws.Range("N67").Value = -33703.5
# Row 102
$ws.Range("H102").Value = 1406.375
$ws.Range("I102").Value = 1223.5555
$ws.Range("J102").Value = 1641.4286
$ws.Range("K102").Value = 1223.5555
$ws.Range("L102").Value = 1641.4286
$ws.Range("M102").Value = 398.4445000000001
$ws.Range("N102").Value = -4885.4286
# Row 126
$ws.Range("H126").Value = 9656.691999999999
$ws.Range("I126").Value = 12892.889
$ws.Range("J126").Value = 2375.25
$ws.Range("K126").Value = 38678.667
$ws.Range("L126").Value = 7125.75
$ws.Range("M126").Value = -36208.667
$ws.Range("N126").Value = -12065.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 48
$ws.Range("H48").Value = 15000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 15000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 15000
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -16322
# Row 54
$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 15000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16288
# Row 122
$ws.Range("H122").Value = 6787404.5
$ws.Range("I122").Value = 10205838
$ws.Range("J122").Value = 2001597
$ws.Range("K122").Value = 30617514
$ws.Range("L122").Value = 6004791
$ws.Range("M122").Value = -30615064
# Row 132
$ws.Range("H132").Value = 13483045
$ws.Range("I132").Value = 17630584
$ws.Range("J132").Value = 3542
$ws.Range("K132").Value = 52891752
$ws.Range("L132").Value = 10626
$ws.Range("M132").Value = -52889222
$ws.Range("N132").Value = -15686

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# Row 126
$ws.Range("H126").Value = 871.0909
$ws.Range("I126").Value = 673.1539
$ws.Range("J126").Value = 1157
$ws.Range("K126").Value = 2019.4617
$ws.Range("L126").Value = 3471
$ws.Range("M126").Value = 450.5382999999999
$ws.Range("N126").Value = -8411
# Row 132
$ws.Range("H132").Value = 1058.2115
$ws.Range("I132").Value = 763.3488
$ws.Range("J132").Value = 2467
$ws.Range("K132").Value = 2290.0464
$ws.Range("L132").Value = 7401
$ws.Range("M132").Value = 239.9535999999998
$ws.Range("N132").Value = -12461
# Row 136
$ws.Range("H136").Value = 1102.3
$ws.Range("I136").Value = 559.48
$ws.Range("J136").Value = 2007
$ws.Range("K136").Value = 1678.44
$ws.Range("L136").Value = 6021
$ws.Range("M136").Value = 871.5599999999999
$ws.Range("N136").Value = -11121

